$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G3" = 1.75
    "H3" = 3.3
    "I3" = 5.5
    "J3" = 2.4
    "L3" = 5.5
    "M3" = 1.1
    "N3" = 7
    "U3" = 2.1
    "V3" = 1.67
    "Z3" = 13
    "AC3" = 7
    "AE3" = 19
    "AF3" = 67
    "AJ3" = 51
    "AK3" = 41
    "AO3" = 9.5
    "AQ3" = 34
    "AW3" = 6.5
    "AX3" = 29
    "AZ3" = 101
    "G4" = 2.35
    "I4" = 3.6
    "U4" = 2.38
    "V4" = 1.53
    "W4" = 5.5
    "X4" = 9.5
    "Z4" = 23
    "AD4" = 5.5
    "AG4" = 7
    "AH4" = 15
    "AV4" = 81
    "AW4" = 5
    "M5" = 1.07
    "N5" = 9
    "G6" = 2.6
    "I6" = 3
    "AX6" = 19
    "G7" = 2.2
    "H7" = 2.88
    "I7" = 3.9
    "S7" = 1.73
    "T7" = 2.08
    "U7" = 2.63
    "V7" = 1.44
    "AB7" = 51
    "AC7" = 5
    "AO7" = 15
    "AU7" = 11
    "G11" = 4.45
    "J11" = 4.8
    "L11" = 2.4
    "S11" = 1.44
    "U11" = 2.02
    "V11" = 1.62
    "X11" = 23
    "Y11" = 15
    "AA11" = 50
    "AC11" = 7.5
    "AD11" = 6.5
    "AG11" = 5.5
    "AH11" = 7.4
    "AI11" = 8.75
    "AJ11" = 14
    "AL11" = 37
    "AO11" = 26
    "AQ11" = 150
    "AR11" = 200
    "AT11" = 2.37
    "BA11" = 75
    "I25" = 4
    "K25" = 2.2
    "O25" = 1.29
    "P25" = 3.5
    "Q25" = 2
    "R25" = 1.85
    "AA25" = 15
    "BA25" = 81
    "H26" = 6
    "I26" = 6.5
    "N26" = 26
    "Q26" = 1.3
    "R26" = 3.5
    "AG26" = 29
    "AP26" = 12
    "G28" = 2.22
    "I28" = 3.05
    "J28" = 2.82
    "K28" = 2.07
    "L28" = 3.65
    "M28" = 1.08
    "N28" = 6.7
    "O28" = 1.35
    "P28" = 2.95
    "T28" = 2.7
    "W28" = 7.4
    "X28" = 10.75
    "Y28" = 9
    "Z28" = 23
    "AA28" = 18.5
    "AB28" = 30
    "AC28" = 6.7
    "AF28" = 70
    "AH28" = 15.5
    "AI28" = 11
    "AJ28" = 40
    "AK28" = 28
    "AL28" = 37
    "AM28" = 600
    "AN28" = 4.2
    "AO28" = 11.75
    "AP28" = 19.5
    "AQ28" = 45
    "AR28" = 75
    "AT28" = 2.7
    "AU28" = 6.9
    "AV28" = 60
    "AW28" = 5
    "AX28" = 17
    "AY28" = 24
    "AZ28" = 80
    "BA28" = 120
    "BB28" = 300
    "H29" = 3.45
    "N29" = 7.9
    "O29" = 1.26
    "AC29" = 7.9
    "AL29" = 25
    "AO29" = 14.5
    "AU29" = 6.9
    "AV29" = 55
    "H30" = 3.75
    "I30" = 3.3
    "L30" = 3.7
    "O30" = 1.16
    "P30" = 4.55
    "Q30" = 1.5
    "R30" = 2.4
    "S30" = 1.28
    "U30" = 1.47
    "V30" = 2.5
    "W30" = 11
    "X30" = 11.75
    "AA30" = 13
    "AD30" = 7.9
    "AG30" = 15
    "AH30" = 22
    "AI30" = 11.75
    "AL30" = 26
    "AS30" = 120
    "AW30" = 5.7
    "AY30" = 19.5
    "BA30" = 90
    "H31" = 3.65
    "I31" = 4.25
    "J31" = 2.3
    "L31" = 4.7
    "N31" = 7.5
    "S31" = 1.4
    "U31" = 1.85
    "V31" = 1.85
    "X31" = 7.9
    "Z31" = 13
    "AC31" = 7.5
    "AD31" = 7.1
    "AE31" = 16.5
    "AF31" = 80
    "AJ31" = 70
    "AP31" = 18.5
    "AQ31" = 29
    "AR31" = 65
    "AU31" = 7.8
    "AV31" = 75
    "AW31" = 6.1
    "I32" = 2.87
    "J32" = 2.85
    "P32" = 3.3
    "W32" = 8.5
    "X32" = 11.75
    "AB32" = 26
    "AG32" = 9.75
    "AH32" = 15.5
    "AJ32" = 35
    "AL32" = 30
    "AP32" = 19
    "AQ32" = 45
    "AR32" = 75
    "AU32" = 6.9
    "AV32" = 60
    "AX32" = 16
    "AZ32" = 75
    "L34" = 3.25
    "N34" = 6.9
    "Q34" = 2
    "R34" = 1.75
    "U34" = 1.78
    "V34" = 1.93
    "W34" = 8
    "AB34" = 30
    "AC34" = 6.9
    "AD34" = 6.2
    "AG34" = 8.25
    "AK34" = 23
    "AP34" = 20
    "AU34" = 6.9
    "AX34" = 14.5
    "AY34" = 22
    "AZ34" = 65
    "BA34" = 100
    "G36" = 2.67
    "H36" = 2.95
    "I36" = 2.67
    "J36" = 3.3
    "L36" = 3.2
    "M36" = 1.03
    "N36" = 12
    "O36" = 1.28
    "P36" = 3.05
    "Q36" = 1.88
    "R36" = 1.82
    "W36" = 8.75
    "X36" = 14
    "Y36" = 9.5
    "Z36" = 32
    "AA36" = 22
    "AB36" = 28
    "AC36" = 9.25
    "AD36" = 5.8
    "AG36" = 9.25
    "AH36" = 14.5
    "AI36" = 9.5
    "AJ36" = 32
    "AN36" = 4.6
    "AO36" = 15
    "AP36" = 21
    "AQ36" = 70
    "AW36" = 4.6
    "AX36" = 14.5
    "AY36" = 20
    "AZ36" = 65
    "M37" = 1.02
    "N37" = 7.1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value2 = $updates[$addr]
}
